# Atualiza os valores 'por municipio' da planilha 'Valores' (relatorio_neomater_APENAS_VALORES)
# criando/preenchendo a matriz A1:L16 com os novos totais calculados.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores")

$data = @(
    @(1,4,1,5,0,2,0,0,3,0,7,0),
    @(0,0,0,0,1,1,0,3,0,5,1,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,1,0,0,0,0,0,0,1,0),
    @(0,3,0,0,0,0,0,0,0,0,0,0),
    @(1,5,1,4,0,2,0,0,4,0,6,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,1,0,1,0,0),
    @(0,0,0,0,2,1,0,1,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,1,0,0,0,0,1,1,1),
    @(2,12,2,10,4,6,0,5,7,7,16,2)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}
